$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings (e.g. "23.551.58", "1.001") that
# must stay literal text rather than being auto-coerced to numbers by Excel.
$dCells = @("D2","D3","D4","D6","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D20","D21","D22","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "23.551.58"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.651.37"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "300.39"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").Value = "0.3788"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").Value = "50.81"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").Value = "0.3567"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "0.08116"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "1.225"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "22.08"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "6.411"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "7.404"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "1.654.32"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "97.17"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "0.06992"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "6.807"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "17.49"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "23.561.30"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").Value = "2.933"
$ws.Range("E26").Value = "  -6.10%  "
$ws.Range("D27").Value = "21.02"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "152.86"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "5.239"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "133.05"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "1.836.79"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").Value = "6.970"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").Value = "2.155"
$ws.Range("E33").Value = "  +5.81%  "
$ws.Range("D34").Value = "11.78"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").Value = "1.041"
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("D36").Value = "0.02745"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").Value = "0.08746"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "0.2451"
$ws.Range("D39").Value = "5.988"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "13.16"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("D41").Value = "0.06891"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").Value = "0.6934"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").Value = "1.319"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "15.73"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "0.6452"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "2.270"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").Value = "3.931"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "0.07872"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "126.90"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  +0.00%  "
